$d = $word.ActiveDocument

# 1. Insert "Test " before "development." in the Python/JS/Mobile skills bullet.
$d.Content.Find.Execute("5+ years Experience in Python, JS, and Mobile development.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "5+ years Experience in Python, JS, and Mobile Test development.", 2) | Out-Null

# 2. Insert "and Mobile app Testing" after "testing Modem/Router firmware" in the Minim summary.
$d.Content.Find.Execute("in the realm of testing Modem/Router firmware. This role allowed me", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "in the realm of testing Modem/Router firmware and Mobile app Testing. This role allowed me", 2) | Out-Null
